# Word COM-interop script: update "Friday July 12, 2013" entry and append the
# "Saturday July 13, 2013" entry to the Weekly Summary document.
#
# The target paragraph is the last paragraph in the document body, so its
# Range does not include the final (un-deletable) paragraph mark. We locate
# that paragraph, then replace its contents with freshly-built OOXML via
# Range.InsertXML -- this lets us reproduce the exact run layout (several
# sibling <w:r> runs with identical formatting, plus a new blank paragraph
# and the "Saturday" paragraph with the trailing bookmark) that the diff
# calls for, rather than relying on Word's automatic run-merging behaviour.

$d = $word.ActiveDocument

# Locate the paragraph that starts the "Friday July 12, 2013" entry -- it is
# the final paragraph of the document.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Friday July 12, 2013*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Friday July 12, 2013' paragraph"
}

$r = $d.Range($target.Range.Start, $target.Range.End)

$xmlSnippet = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00E11095" w:rsidRPr="00C9751A" w:rsidRDefault="00E11095" w:rsidP="00C9751A">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Friday July 12, 2013 </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>–</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Took a day off to focus on other classes.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Saturday July 13, 2013 – Figured out how to alter many pixels at a time or a single pixel. Found out in my research that each pixel found with x and y coordinates has bands. Bands are consisted of three different values that are red, green, and blue. There can be more or less than three bands, but almost every image consists of only the RGB values.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> Scroll bars are now usable. </w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$r.InsertXML($xmlSnippet)
